# Applies the "Story for the areas" edit: expands the one-line area
# teaser paragraphs into their full intro-run + detail-run(s) form, and
# splits two existing runs (Desert / Cave areas) into several runs.
#
# Strategy: Word's COM `Range.InsertAfter` live-edits the run model and
# silently re-merges adjacent runs that share formatting, so sequential
# InsertAfter calls can't reproduce the multi-run split seen in the
# target XML. `Range.InsertXML`, however, *replaces* the exact range it
# is called on with literal OOXML content (preserving the surrounding
# paragraph's own attributes) and does not re-merge the runs it creates.
# So every edit below is expressed as: find the exact sub-range of text
# to change, then call InsertXML with a small single-paragraph package
# whose runs become the new run sequence at that spot.

$d = $word.ActiveDocument

function Escape-Xml([string]$s) {
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

function Run-Xml([string]$text) {
    $needsSpace = ($text.Length -eq 0) -or ($text[0] -eq ' ') -or ($text[$text.Length - 1] -eq ' ')
    $escaped = Escape-Xml $text
    if ($needsSpace) {
        return "<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
    } else {
        return "<w:r><w:t>$escaped</w:t></w:r>"
    }
}

function Package-Xml([string]$innerRunsXml) {
    return '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        "<w:body><w:p>$innerRunsXml</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
}

# Replace the literal substring $oldText (which must occur exactly once
# inside paragraph $paraIndex) with a sequence of freshly split runs
# whose texts are given in $newTexts.
function Replace-TextWithRuns([int]$paraIndex, [string]$oldText, [string[]]$newTexts) {
    $p = $d.Paragraphs.Item($paraIndex)
    $full = $p.Range.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        throw "Could not find target text in paragraph $paraIndex : $oldText"
    }
    $start = $p.Range.Start + $idx
    $end = $start + $oldText.Length
    $target = $d.Range($start, $end)

    $runsXml = ""
    foreach ($t in $newTexts) {
        $runsXml += Run-Xml $t
    }
    $target.InsertXML((Package-Xml $runsXml))
}

# Desert area: split the second run into three runs.
Replace-TextWithRuns 6 `
    "An enormous sandworm has an object inside of it which will help to create the ultimate weapon in the late game." `
    @(
        "An enormous sandworm has ",
        "a part of a powerful sword ",
        "inside of it which will help to create the ultimate weapon in the late game."
    )

# Cave area: split the second run into three runs.
Replace-TextWithRuns 7 `
    ": Rumored there is a potion in the caves which can heal any curse or wound. But there is also its guardian protecting it from harm." `
    @(
        ": Rumored there is a potion in the caves which can heal any curse or wound. But there is also its guardian protecting it from ",
        "thieves",
        "."
    )

# Cold area: append one new run with the full sentence.
Replace-TextWithRuns 8 "Cold area" `
    @(
        "Cold area",
        ": Legends say there is a blacksmith in the cold north who can fix every weapon possible if you can do his challenges."
    )

# Lake area: append one new run with the full sentence.
Replace-TextWithRuns 9 "Lake area" `
    @(
        "Lake area",
        ": Sunken in the largest lake of the land there is a hidden metal which is more powerful than everything else in the world."
    )

# Ancient area: append two new runs.
Replace-TextWithRuns 10 "Ancient area" `
    @(
        "Ancient area",
        ": It is handed down that there is an ancient handle in the ruins of the ",
        "elders having the power to hold together the powerfullest of all weapons."
    )

# Volcano area: append two new runs.
Replace-TextWithRuns 11 "Volcano area" `
    @(
        "Volcano area",
        ": ",
        "A special lava in the deepest point of the volcano is supposed to kill the witch who was seemingly invincible."
    )

# Sea area: append four new runs.
Replace-TextWithRuns 12 "Sea area" `
    @(
        "Sea area",
        ": ",
        "An aged gemston",
        "e is protected by an old man living at the sea. He is not only the ",
        "commander of all enemies but also a member of the ancient people"
    )

# Mountain area: append two new runs.
Replace-TextWithRuns 13 "Mountain area" `
    @(
        "Mountain area",
        ": ",
        "The Witches oldest allies are reawakened and must be slain in the mountains before they can reach and support the witch in the last fight."
    )

# Dark forest area: append two new runs.
Replace-TextWithRuns 14 "Dark forest area" `
    @(
        "Dark forest area",
        ": ",
        "The last fight versus the witch is held and after she is slain the swordsman is free and turns back to the village."
    )
